$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 2015
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 171623.3273763794
$ws.Range("D33").Value = 2447.277127511424
$ws.Range("E33").Value = 81.68503211498714
$ws.Range("F33").Value = 221.4275749705685
$ws.Range("G33").Value = 530

$ws.Range("A34").Value = 2015
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 493755.4590511023
$ws.Range("D34").Value = 7040.747082544511
$ws.Range("E34").Value = 235.0055271983661
$ws.Range("F34").Value = 637.0408708274134
$ws.Range("G34").Value = 1070

$ws.Range("A35").Value = 2015
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 510048.4200640724
$ws.Range("D35").Value = 7273.077916815734
$ws.Range("E35").Value = 242.760248331363
$ws.Range("F35").Value = 658.0619691904119
$ws.Range("G35").Value = 1260

$ws.Range("A36").Value = 2015
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 745049.2630870403
$ws.Range("D36").Value = 10624.09200602696
$ws.Range("E36").Value = 354.6101448624124
$ws.Range("F36").Value = 961.2589039083937
$ws.Range("G36").Value = 2030

$ws.Range("A37").Value = 2015
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 235000.843029632
$ws.Range("D37").Value = 3351.014089319564
$ws.Range("E37").Value = 111.8498965347412
$ws.Range("F37").Value = 303.1969347264659
$ws.Range("G37").Value = 770

$ws.Range("A38").Value = 2015
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 235000.843029632
$ws.Range("D38").Value = 3351.014089319564
$ws.Range("E38").Value = 111.8498965347412
$ws.Range("F38").Value = 303.1969347264659
$ws.Range("G38").Value = 770

$ws.Range("A39").Value = 2015
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 235000.843029632
$ws.Range("D39").Value = 3351.014089319564
$ws.Range("E39").Value = 111.8498965347412
$ws.Range("F39").Value = 303.1969347264659
$ws.Range("G39").Value = 770

$ws.Range("A40").Value = 2015
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 235000.843029632
$ws.Range("D40").Value = 3351.014089319564
$ws.Range("E40").Value = 111.8498965347412
$ws.Range("F40").Value = 303.1969347264659
$ws.Range("G40").Value = 770
